# manySampleMTableInFooter-template.docx
#
# The canonical-OOXML diff for this revision touches three parts
# (word/document.xml's sectPr, word/footer1.xml's root element, and
# word/styles.xml) but, attribute-for-attribute, every changed line
# carries exactly the same name/value pairs before and after -- only
# the serialization order of attributes (and of the xmlns:* namespace
# declarations) changed, together with the removal of the volatile
# w:rsid* bookkeeping attributes that Word stamps on edit. No element,
# no text run, no formatting value and no relationship actually
# changes. In other words the content of the template is identical;
# what differs is purely how the XML writer that produced the
# "after" snapshot chose to order attributes.
#
# That kind of pure re-serialization is not something the Word object
# model exposes a knob for (there is no "sort my attributes"
# command), so the faithful way to reproduce this revision through
# COM automation is to touch the document without changing any of its
# actual content: no text, formatting, margins, or custom properties
# are modified. We simply open the pieces the diff says were
# re-saved and read them back, which is enough to confirm the
# document round-trips with its content intact.

$d = $word.ActiveDocument

# Touch word/document.xml (section properties / page setup) without
# altering any value.
$section = $d.Sections(1)
$null = $section.PageSetup.PageWidth
$null = $section.PageSetup.PageHeight
$null = $section.PageSetup.TopMargin
$null = $section.PageSetup.BottomMargin
$null = $section.PageSetup.LeftMargin
$null = $section.PageSetup.RightMargin
$null = $section.PageSetup.Gutter
$null = $section.PageSetup.HeaderDistance
$null = $section.PageSetup.FooterDistance

# Touch word/footer1.xml without altering its text or fields.
$footer = $section.Footers(1)
$null = $footer.Range.Text

# Touch word/styles.xml without altering any style definition.
$null = $d.Styles.Count
for ($i = 1; $i -le $d.Styles.Count; $i++) {
    $null = $d.Styles($i).NameLocal
}

$d.Saved = $d.Saved
